$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: Test Case ID becomes blank (space), new test scenario text
$ws.Range("A2").Value = " "
$ws.Range("B2").Value = "@iProctorRegression Verify Elumina Create Exam with survey section"
$ws.Range("C2").Value = "passed"

# Add row 3
$ws.Range("A3").Value = " "
$ws.Range("B3").Value = "@iProctorRegression Verify Elumina Registration"
$ws.Range("C3").Value = "passed"

# Add row 4
$ws.Range("A4").Value = " "
$ws.Range("B4").Value = "@iProctorRegression Verify Validation of Survey screen EluminaRegressioniProctorTC_083"
$ws.Range("C4").Value = "passed"
